$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "newexam" package.
$ws.Name = "newexam.cases.TestNewExamPaper"

# Move the active selection to I15 (was C6).
$ws.Range("I15").Select()
